# Weekly data refresh: a new daily price record is inserted as row 75
# (pushing the existing rows 75-151 down to 76-152), matching the
# "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75; Excel shifts rows 75..151 down to 76..152
# and the sheet dimension grows from R151 to R152 automatically.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44587
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112024
$ws.Range("G75").Value = "Choclo"
$ws.Range("H75").Value = "Choclero"
$ws.Range("I75").Value = "Segunda"
$ws.Range("J75").Value = 1000
$ws.Range("K75").Value = 150
$ws.Range("L75").Value = 200
$ws.Range("M75").Value = 175
$ws.Range("N75").Value = "$/unidad"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 175
$ws.Range("Q75").Value = 1
$ws.Range("R75").Value = "Hortaliza"
